# Fixed Bento 80 Test scripts
#
# Each row of the "startup" sheet holds a Cypher query (column B) used by
# one of the Bento tabs (CasesTab / SamplesTab / FilesTab, named in column
# A). This appends a deterministic "ORDER BY ... LIMIT 100" clause to each
# query's text, and updates the saved selection from B4 to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$orderClauses = @{
    "FilesTab"   = "`n order By f.file_name ASC LIMIT 100"
    "SamplesTab" = "`n order By samp.sample_id ASC LIMIT 100"
    "CasesTab"   = "`n order By ss.study_subject_id ASC LIMIT 100 "
}

# Apply to FilesTab, then SamplesTab, then CasesTab (in that order) so the
# workbook's shared-string table is rebuilt in the same order Excel wrote it.
foreach ($tabName in @("FilesTab", "SamplesTab", "CasesTab")) {
    for ($row = 1; $row -le 4; $row++) {
        if ($ws.Cells.Item($row, 1).Value() -eq $tabName) {
            $query = $ws.Cells.Item($row, 2).Value()
            $ws.Cells.Item($row, 2).Value = $query + $orderClauses[$tabName]
            break
        }
    }
}

# Move the active selection on the sheet from B4 to B2.
$ws.Range("B2").Select()
